# Final user list changes - 16 May 2025
#
# Users!A2 is updated from "Christine Goodridge" to "Amanda Donovan".
# The previously-selected/active sheet (UpdateActivity) is replaced by the
# Users sheet becoming the active tab, with its selection moved to E8.

$wb = $excel.ActiveWorkbook

# Update the user name on the Users sheet.
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Amanda Donovan"

# Make "Users" the active sheet/tab, with E8 selected - mirrors the
# workbook now opening on the Users sheet instead of UpdateActivity.
$wsUsers.Activate()
$wsUsers.Range("E8").Select()
